$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- No. 1 (multiples of 3 or 5) -----------------------------------------
$ws.Range("B3").Value = "10보다 작은 자연수 중에서 3 또는 5의 배수는 3, 5, 6, 9 이고, 이것을 모두 더하면 23입니다.`n1000보다 작은 자연수 중에서 3 또는 5의 배수를 모두 더하면 얼마일까요?"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "1000보다 작은수니까 i<=1000 이 아니라 i<1000 이지…"
[void]$ws.Rows.Item(3).AutoFit()

# --- No. 2 (even Fibonacci numbers) --------------------------------------
$ws.Range("B4").Value = "피보나치 수열의 각 항은 바로 앞의 항 두 개를 더한 것이 됩니다. 1과 2로 시작하는 경우 이 수열은 아래와 같습니다.`n1, 2, 3, 5, 8, 13, 21, 34, 55, 89, ...`n짝수이면서 4백만 이하인 모든 항을 더하면 얼마가 됩니까?"
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = "While조건의 꺽쇠 방향…"
[void]$ws.Rows.Item(4).AutoFit()

# --- Updated times for the already-solved problems (rows 5 & 6) ---------
$ws.Range("C5").Value = 35
$ws.Range("C6").Value = 30

# Row 5 shrinks now that its feedback note reflows into fewer lines
$ws.Rows.Item(5).RowHeight = 33

# --- Selection / scroll position left by the author after editing -------
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
[void]$ws.Range("D13").Select()
